$d = $word.ActiveDocument

# --- Step 1: Split "Miska Sainkangas" into two runs with a spell-check
#     proofErr wrapper around "Sainkangas" (leave paragraph's pPr/rsid attrs
#     untouched by only replacing the inner run text, not the whole <w:p>) ---

$miskaPara = $d.Paragraphs.Item(6)
$miskaStart = $miskaPara.Range.Start
$miskaTextRange = $d.Range($miskaStart, $miskaStart + 16)  # "Miska Sainkangas" (no para mark)

$miskaFrag = '<w:p><w:r><w:t xml:space="preserve">Miska </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Sainkangas</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$miskaPkg = "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body>$miskaFrag</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
$miskaTextRange.InsertXML($miskaPkg)

# Remove the old (now orphaned-in-place) "_GoBack" bookmark from this paragraph;
# it will be re-created further down, at the new edit location.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- Step 2: After the page-break run, add a new run of text, preserving
#     the lastRenderedPageBreak marks on both runs, then drop a fresh
#     "_GoBack" bookmark right after the inserted text. ---

$pageBreakPara = $d.Paragraphs.Item(20)
$pbStart = $pageBreakPara.Range.Start
$pbEnd = $pageBreakPara.Range.End
$pageBreakCharRange = $d.Range($pbStart, $pbEnd - 1)  # the page-break char only, excludes the paragraph mark

$pbFrag = '<w:p><w:r><w:lastRenderedPageBreak/><w:br w:type="page"/></w:r><w:r><w:lastRenderedPageBreak/><w:t>gpoijdsgoedsgoidsjfrgoidjgroidrgoijd</w:t></w:r></w:p>'
$pbPkg = "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body>$pbFrag</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
$pageBreakCharRange.InsertXML($pbPkg)

# Re-fetch the paragraph/position, since the document shifted after the insert.
$pageBreakPara = $d.Paragraphs.Item(20)
$newBookmarkPos = $pageBreakPara.Range.End - 1
$newBookmarkRange = $d.Range($newBookmarkPos, $newBookmarkPos)
$d.Bookmarks.Add("_GoBack", $newBookmarkRange)
